# Update the "update-employee" format worksheet:
#  - rename two existing header labels (Code -> Name variants)
#  - insert a new "Branch" column before the trailing "Aadhaar" column
#  - adjust the new/old last-column widths
#  - restore the active selection on the newly added column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename Designation Code / Department Code headers to the *Name variants
$ws.Range("L1").Value = "Designation Name"
$ws.Range("M1").Value = "Department Name"

# 2) Insert a new column before the last one (Aadhaar, currently column AB)
#    so Aadhaar shifts right to AC and the new column becomes AB.
$ws.Columns.Item(28).Insert()
$ws.Range("AB1").Value = "Branch"

# 3) Column width bookkeeping: the new AB column keeps the wider (21.625)
#    formatting while the shifted Aadhaar column (now AC) keeps its old
#    narrower numeric-style width (18.625). (ColumnWidth is expressed in
#    Excel's character-width units, ~5/6 less than the raw sheet width.)
$ws.Columns.Item(28).ColumnWidth = 20.791666666666668
$ws.Columns.Item(29).ColumnWidth = 17.791666666666668

# 4) Move the active selection to the newly inserted column header cell.
$ws.Range("AB1").Select() | Out-Null
